$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this shifts existing rows 28:100 down to 29:101
# (carrying their values/styles with them), matching the target diff where every
# row from 29..101 (after) equals row 28..100 (before), and the dimension grows
# from A1:R100 to A1:R101.
$ws.Rows("28").Insert()

# Populate the newly inserted row 28 with the new record.
$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 45177
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 300000000
$ws.Range("G28").Value = "Espárragos"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 2500
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2500
$ws.Range("N28").Value = "$/kilo"
$ws.Range("O28").Value = "Provincia de Linares"
$ws.Range("P28").Value = 2500
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"
